$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "cxq6hz_20170224_144343_ASIC_EEG"
$ws.Range("F1").Value = "hyy-调节6Hz_20170306_110203_ASIC_EEG"
$ws.Range("G1").Value = "hzj-调节6Hz_20170220_113105_ASIC_EEG"

$ws.Range("E2").Value = 0.63501483679525217
$ws.Range("F2").Value = 0.64646464646464641
$ws.Range("G2").Value = 0.62462462462462465

$ws.Range("E3").Value = 0.65597667638483959
$ws.Range("F3").Value = 0.65187713310580198
$ws.Range("G3").Value = 0.65161290322580645

$ws.Range("A1:G3").Select()

$wb.Save()
